$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "93.388.43"
$ws.Range("E2").Value = "  +3.04%  "
$ws.Range("D3").Value = "3.129.34"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "616.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("E7").Value = "  +1.64%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.407"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +10.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "3.130.14"
$ws.Range("E10").Value = "  +31.03%  "
$ws.Range("E11").Value = "  +1.58%  "
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("E13").Value = "  +4.33%  "
$ws.Range("D14").Value = "93.239.43"
$ws.Range("E14").Value = "  +3.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "3.727.45"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").Value = "3.148.80"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("E19").Value = "  +4.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("E22").Value = "  +3.54%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "450.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.40%  "
$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.42%  "
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "87.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("D28").Value = "3.298.94"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  +13.61%  "
$ws.Range("E31").Value = "  +1.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.170"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("E35").Value = "  +7.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.160"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.58%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "497.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.20%  "
$ws.Range("B39").Value = "PancakeSwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("E40").Value = "  +5.31%  "
$ws.Range("E41").Value = "  -2.47%  "
$ws.Range("E42").Value = "  -1.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E44").Value = "  +4.51%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.19%  "
$ws.Range("E47").Value = "  +1.56%  "
$ws.Range("E48").Value = "  -1.62%  "
$ws.Range("E49").Value = "  +3.75%  "
$ws.Range("E50").Value = "  +5.20%  "
$ws.Range("E51").Value = "  +1.13%  "
